# Update "paises" (countries) COVID data sheet + refresh timestamp.
# Matches commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the "last updated" timestamp label in A1.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 01:52"

# ---------------------------------------------------------------------
# Helper-less direct cell writes: row layout is
#   A = Pais, B = Casos totales, C = Nuevos casos, D = Casos activos,
#   E = Recuperados, F = Casos criticos, G = Muertes hoy, H = Muertes
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 2) Countries whose stats changed but whose sorted rank stays the same.
# ---------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Cells.Item(4,2).Value = 1063351
$ws.Cells.Item(4,3).Value = 27586
$ws.Cells.Item(4,4).Value = 147114
$ws.Cells.Item(4,5).Value = 854619
$ws.Cells.Item(4,6).Value = 18671
$ws.Cells.Item(4,7).Value = 2352
$ws.Cells.Item(4,8).Value = 61618

# Row 14 - Brasil
$ws.Cells.Item(14,2).Value = 79361
$ws.Cells.Item(14,3).Value = 6462
$ws.Cells.Item(14,4).Value = 34132
$ws.Cells.Item(14,5).Value = 39579
$ws.Cells.Item(14,6).Value = 8318
$ws.Cells.Item(14,7).Value = 448
$ws.Cells.Item(14,8).Value = 5511

# Row 15 - Canada
$ws.Cells.Item(15,2).Value = 51597
$ws.Cells.Item(15,3).Value = 1571
$ws.Cells.Item(15,4).Value = 20327
$ws.Cells.Item(15,5).Value = 28274
$ws.Cells.Item(15,6).Value = 557
$ws.Cells.Item(15,7).Value = 137
$ws.Cells.Item(15,8).Value = 2996

# Row 44 - Noruega
$ws.Cells.Item(44,2).Value = 7710
$ws.Cells.Item(44,3).Value = 50
$ws.Cells.Item(44,4).Value = 32
$ws.Cells.Item(44,5).Value = 7471
$ws.Cells.Item(44,6).Value = 40
$ws.Cells.Item(44,7).Value = 1
$ws.Cells.Item(44,8).Value = 207

# Row 56 - Argentina
$ws.Cells.Item(56,2).Value = 4285
$ws.Cells.Item(56,3).Value = 158
$ws.Cells.Item(56,4).Value = 1192
$ws.Cells.Item(56,5).Value = 2879
$ws.Cells.Item(56,6).Value = 144
$ws.Cells.Item(56,7).Value = 7
$ws.Cells.Item(56,8).Value = 214

# Row 87 - Costa de Marfil
$ws.Cells.Item(87,2).Value = 1238
$ws.Cells.Item(87,3).Value = 55
$ws.Cells.Item(87,4).Value = 557
$ws.Cells.Item(87,5).Value = 667
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 0
$ws.Cells.Item(87,8).Value = 14

# Row 128 - Maldivas
$ws.Cells.Item(128,2).Value = 278
$ws.Cells.Item(128,3).Value = 28
$ws.Cells.Item(128,4).Value = 17
$ws.Cells.Item(128,5).Value = 260
$ws.Cells.Item(128,6).Value = 2
$ws.Cells.Item(128,7).Value = 1
$ws.Cells.Item(128,8).Value = 1

# ---------------------------------------------------------------------
# 3) Nigeria's total-case count overtakes Bosnia/Ghana/Estonia, so it
#    moves up to just before "Bosnia y Herzegovina" (old row 75).
#    Insert its new row there, then remove its old row (now pushed one
#    row further down, at 79, since the insert shifted everything by 1).
# ---------------------------------------------------------------------
$ws.Rows.Item(75).Insert()
$ws.Cells.Item(75,1).Value = "Nigeria"
$ws.Cells.Item(75,2).Value = 1728
$ws.Cells.Item(75,3).Value = 196
$ws.Cells.Item(75,4).Value = 307
$ws.Cells.Item(75,5).Value = 1370
$ws.Cells.Item(75,6).Value = 2
$ws.Cells.Item(75,7).Value = 7
$ws.Cells.Item(75,8).Value = 51
$ws.Rows.Item(79).Delete()

# ---------------------------------------------------------------------
# 4) San Vicente y las Granadinas overtakes Dominica/Curazao/San
#    Cristobal y Nieves, so it moves up to just before "Dominica"
#    (old row 192). Same insert+delete dance as above.
# ---------------------------------------------------------------------
$ws.Rows.Item(192).Insert()
$ws.Cells.Item(192,1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(192,2).Value = 16
$ws.Cells.Item(192,3).Value = 1
$ws.Cells.Item(192,4).Value = 8
$ws.Cells.Item(192,5).Value = 8
$ws.Cells.Item(192,6).Value = 0
$ws.Cells.Item(192,7).Value = 0
$ws.Cells.Item(192,8).Value = 0
$ws.Rows.Item(196).Delete()
